$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E: plain text updates (never numeric-looking, safe as-is)
$plainUpdates = @(
    @('E3', '  +1.23%  '),
    @('E4', '  +0.40%  '),
    @('E5', '  +0.35%  '),
    @('E6', '  +0.10%  '),
    @('E7', '  +2.57%  '),
    @('E8', '  +2.38%  '),
    @('E9', '  +1.31%  '),
    @('E10', '  +3.16%  '),
    @('E11', '  +1.63%  '),
    @('E12', '  +0.26%  '),
    @('B13', 'WrappedEther'),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('E13', '  +4.60%  '),
    @('B14', 'Polkadot'),
    @('C14', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'),
    @('E14', '  +1.37%  '),
    @('E15', '  +1.21%  '),
    @('E16', '  +1.47%  '),
    @('E17', '  +0.47%  '),
    @('E18', '  +0.98%  '),
    @('E19', '  +0.39%  '),
    @('E20', '  +1.34%  '),
    @('E21', '  +2.10%  '),
    @('E22', '  +0.63%  '),
    @('E23', '  +0.94%  '),
    @('E24', '  +1.30%  '),
    @('E25', '  -0.18%  '),
    @('E26', '  +1.45%  '),
    @('E27', '  +0.91%  '),
    @('E28', '  +0.40%  '),
    @('E29', '  +2.75%  '),
    @('E30', '  +1.04%  '),
    @('E31', '  +5.33%  '),
    @('E32', '  +0.87%  '),
    @('E33', '  +1.53%  '),
    @('E34', '  +1.62%  '),
    @('E35', '  -3.52%  '),
    @('E36', '  +3.05%  '),
    @('E37', '  +1.33%  '),
    @('E38', '  +2.00%  '),
    @('E39', '  +2.19%  '),
    @('E40', '  +3.28%  '),
    @('E41', '  +1.83%  '),
    @('E42', '  +1.32%  '),
    @('E43', '  +1.78%  '),
    @('E44', '  +5.71%  '),
    @('E45', '  +2.01%  '),
    @('E46', '  +0.41%  '),
    @('E47', '  +4.09%  '),
    @('E48', '  +2.62%  '),
    @('E49', '  +2.53%  '),
    @('E50', '  +0.68%  '),
    @('E51', '  +4.09%  ')
)
foreach ($u in $plainUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Column D: values that look numeric (e.g. "1.011", "309.92") must be forced
# to text so Excel does not reinterpret them as numbers, matching the
# original inline-string / shared-string text cell type. We temporarily set
# a text NumberFormat, assign the value, then clear the format again so the
# cell keeps its original (unstyled) appearance.
$textUpdates = @(
    @('D2', '26.995.43'),
    @('D3', '1.848.32'),
    @('D5', '1.011'),
    @('D6', '309.92'),
    @('D7', '0.4783'),
    @('D8', '0.3683'),
    @('D9', '0.07235'),
    @('D10', '0.9317'),
    @('D11', '19.75'),
    @('D12', '0.07718'),
    @('D13', '1.907.97'),
    @('D14', '5.341'),
    @('D15', '6.441'),
    @('D17', '1.014'),
    @('D18', '0.000008643'),
    @('D20', '27.028.72'),
    @('D21', '14.50'),
    @('D22', '5.057'),
    @('D24', '1.930'),
    @('D25', '152.79'),
    @('D27', '2.004'),
    @('D28', '114.40'),
    @('D29', '5.003'),
    @('D30', '0.08904'),
    @('D31', '3.292'),
    @('D32', '1.177'),
    @('D33', '0.7452'),
    @('D34', '4.510'),
    @('D35', '2.743'),
    @('D36', '1.115'),
    @('D37', '0.01957'),
    @('D38', '0.05266'),
    @('D40', '0.5229'),
    @('D41', '7.003'),
    @('D42', '0.1517'),
    @('D43', '8.218'),
    @('D44', '10.63'),
    @('D45', '0.4756'),
    @('D47', '102.09'),
    @('D48', '1.613'),
    @('D49', '65.55'),
    @('D50', '0.06075'),
    @('D51', '0.8887')
)
foreach ($u in $textUpdates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.ClearFormats()
}
